$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '28.491.79'
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.567.82'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -2.12%  '
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '211.89'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -1.49%  '
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.493'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.17%  '
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.ClearFormats()

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '46.15'
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -1.58%  '
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0887'
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.40%  '
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.791.30'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -2.18%  '
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.572.05'
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -2.04%  '
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -2.64%  '
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -3.05%  '
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '28.488.43'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '62.26'
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '227.82'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -1.95%  '
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.36'
$c.ClearFormats()

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -2.57%  '
$c.ClearFormats()

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0₃0692'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -2.82%  '
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.06%  '
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -5.91%  '
$c.ClearFormats()

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '9.12'
$c.ClearFormats()

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -3.25%  '
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +6.17%  '
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '150.97'
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -1.16%  '
$c.ClearFormats()

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '15.00'
$c.ClearFormats()

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -2.12%  '
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -2.65%  '
$c.ClearFormats()

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -4.14%  '
$c.ClearFormats()

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.ClearFormats()

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -1.77%  '
$c.ClearFormats()

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.10'
$c.ClearFormats()

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -4.18%  '
$c.ClearFormats()

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -1.33%  '
$c.ClearFormats()

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.393.63'
$c.ClearFormats()

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -2.19%  '
$c.ClearFormats()

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '
$c.ClearFormats()

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -3.54%  '
$c.ClearFormats()

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.ClearFormats()

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.ClearFormats()

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +2.17%  '
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -0.96%  '
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.537'
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -1.37%  '
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.ClearFormats()

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.88'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +1.83%  '
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -4.50%  '
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -4.23%  '
$c.ClearFormats()

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.974'
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '62.93'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -3.05%  '
$c.ClearFormats()

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.704.08'
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.16%  '
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '85.92'
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.30%  '
$c.ClearFormats()

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -4.81%  '
$c.ClearFormats()

